$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated roster data (player, position, team) replacing the prior table contents
$data = @(
    @("Chris Paul",          "PG",       "San Antonio Spurs"),
    @("Russell Westbrook",   "PG,SG",    "Denver Nuggets"),
    @("Payton Pritchard",    "PG,SG",    "Boston Celtics"),
    @("Brandin Podziemski",  "SG",       "Golden State Warriors"),
    @("Ayo Dosunmu",         "PG,SG,SF", "Chicago Bulls"),
    @("Pascal Siakam",       "SF,PF,C",  "Indiana Pacers"),
    @("Jaylen Brown",        "SG,SF",    "Boston Celtics"),
    @("Deni Avdija",         "SF,PF",    "Portland Trail Blazers"),
    @("Nikola Jokic",        "C",        "Denver Nuggets"),
    @("Rudy Gobert",         "C",        "Minnesota Timberwolves"),
    @("Dejounte Murray",     "PG,SG",    "New Orleans Pelicans"),
    @("Jakob Poeltl",        "C",        "Toronto Raptors"),
    @("Khris Middleton",     "SF",       "Milwaukee Bucks"),
    @("Paolo Banchero",      "SF,PF",    "Orlando Magic"),
    @("Jalen Green",         "PG,SG",    "Houston Rockets"),
    @("Chet Holmgren",       "PF,C",     "Orlando Magic"),
    @("Jalen Suggs",         "PG,SG",    "Orlando Magic")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
